$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case "de/del/la/las/los/el/y" connector words in state/municipality names
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B7").Value = "San Francisco De Los Romo"
$ws.Range("B11").Value = "Playas De Rosarito"
$ws.Range("B22").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Benemérito De Las Américas"
$ws.Range("B70").Value = "Guadalupe Y Calvo"
$ws.Range("B72").Value = "Hidalgo Del Parral"
$ws.Range("B86").Value = "San Francisco Del Oro"
$ws.Range("B93").Value = "Valle De Zaragoza"
$ws.Range("B114").Value = "Villa De Álvarez"
$ws.Range("A116").Value = "Ciudad De México"
$ws.Range("B120").Value = "Cuajimalpa De Morelos"
$ws.Range("B134").Value = "Coneto De Comonfort"
$ws.Range("B146").Value = "Nombre De Dios"
$ws.Range("B149").Value = "Pánuco De Coronado"
$ws.Range("B154").Value = "San Juan De Guadalupe"
$ws.Range("B155").Value = "San Juan Del Río"
$ws.Range("B156").Value = "San Luis Del Cordero"
$ws.Range("A164").Value = "Estado De México"
$ws.Range("B164").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B167").Value = "Almoloya De Alquisiras"
$ws.Range("B170").Value = "Atizapán De Zaragoza"
$ws.Range("B178").Value = "Coacalco De Berriozábal"
$ws.Range("B182").Value = "Ecatepec De Morelos"
$ws.Range("B185").Value = "Ixtapan De La Sal"
$ws.Range("B194").Value = "Naucalpan De Juárez"
$ws.Range("B199").Value = "San Felipe Del Progreso"
$ws.Range("B207").Value = "Tenango Del Valle"
$ws.Range("B211").Value = "Tlalnepantla De Baz"
$ws.Range("B215").Value = "Valle De Bravo"
$ws.Range("B216").Value = "Villa De Allende"
$ws.Range("B217").Value = "Villa Del Carbón"
$ws.Range("B225").Value = "San Miguel De Allende"
$ws.Range("B226").Value = "Apaseo El Alto"
$ws.Range("B227").Value = "Apaseo El Grande"
$ws.Range("B232").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B236").Value = "Jaral Del Progreso"
$ws.Range("B246").Value = "San Diego De La Unión"
$ws.Range("B248").Value = "San Francisco Del Rincón"
$ws.Range("B249").Value = "San Luis De La Paz"
$ws.Range("B250").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B251").Value = "Silao De La Victoria"
$ws.Range("B255").Value = "Valle De Santiago"
$ws.Range("B260").Value = "Acapulco De Juárez"
$ws.Range("B262").Value = "Ajuchitlán Del Progreso"
$ws.Range("B263").Value = "Alcozauca De Guerrero"
$ws.Range("B266").Value = "Atenango Del Río"
$ws.Range("B267").Value = "Atlamajalcingo Del Monte"
$ws.Range("B269").Value = "Atoyac De Álvarez"
$ws.Range("B270").Value = "Ayutla De Los Libres"
$ws.Range("B272").Value = "Buenavista De Cuéllar"
$ws.Range("B273").Value = "Chilapa De Álvarez"
$ws.Range("B274").Value = "Chilpancingo De Los Bravo"
$ws.Range("B275").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B278").Value = "Coyuca De Benítez"
$ws.Range("B279").Value = "Coyuca De Catalán"
$ws.Range("B282").Value = "Cuetzala Del Progreso"
$ws.Range("B287").Value = "Huitzuco De Los Figueroa"
$ws.Range("B288").Value = "Iguala De La Independencia"
$ws.Range("B290").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B301").Value = "Taxco De Alarcón"
$ws.Range("B303").Value = "Técpan De Galeana"
$ws.Range("B305").Value = "Tepecoacuilco De Trujano"
$ws.Range("B306").Value = "Tixtla De Guerrero"
$ws.Range("B309").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B310").Value = "Tlapa De Comonfort"
$ws.Range("B316").Value = "Agua Blanca De Iturbide"
$ws.Range("B319").Value = "Atotonilco De Tula"
$ws.Range("B322").Value = "Cuautepec De Hinojosa"
$ws.Range("B326").Value = "Huejutla De Reyes"
$ws.Range("B329").Value = "Jacala De Ledezma"
$ws.Range("B333").Value = "Mineral Del Monte"
$ws.Range("B334").Value = "Mixquiahuala De Juárez"
$ws.Range("B335").Value = "Molango De Escamilla"
$ws.Range("B336").Value = "Omitlán De Juárez"
$ws.Range("B337").Value = "Pachuca De Soto"
$ws.Range("B339").Value = "Progreso De Obregón"
$ws.Range("B342").Value = "Tenango De Doria"
$ws.Range("B343").Value = "Tepehuacán De Guerrero"
$ws.Range("B344").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B346").Value = "Tezontepec De Aldama"
$ws.Range("B350").Value = "Tulancingo De Bravo"
$ws.Range("B351").Value = "Villa De Tezontepec"
$ws.Range("B354").Value = "Zacualtipán De Ángeles"
$ws.Range("B355").Value = "Zapotlán De Juárez"
$ws.Range("B359").Value = "Acatlán De Juárez"
$ws.Range("B360").Value = "Ahualulco De Mercado"
$ws.Range("B365").Value = "Atotonilco El Alto"
$ws.Range("B366").Value = "Autlán De Navarro"
$ws.Range("B373").Value = "Encarnación De Díaz"
$ws.Range("B376").Value = "Ixtlahuacán Del Río"
$ws.Range("B382").Value = "Lagos De Moreno"
$ws.Range("B388").Value = "San Juan De Los Lagos"
$ws.Range("B389").Value = "San Juanito De Escobedo"
$ws.Range("B390").Value = "San Miguel El Alto"
$ws.Range("B392").Value = "Tamazula De Gordiano"
$ws.Range("B395").Value = "Teocuitatlán De Corona"
$ws.Range("B396").Value = "Tepatitlán De Morelos"
$ws.Range("B398").Value = "Tlajomulco De Zúñiga"
$ws.Range("B403").Value = "Unión De San Antonio"
$ws.Range("B407").Value = "Yahualica De González Gallo"
$ws.Range("B410").Value = "Zapotlán El Grande"
$ws.Range("B426").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B474").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B497").Value = "Puente De Ixtla"
$ws.Range("B500").Value = "Tlaltizapán De Zapata"
$ws.Range("B508").Value = "Santa María Del Oro"
$ws.Range("B516").Value = "Ciénega De Flores"
$ws.Range("B520").Value = "Mier Y Noriega"
$ws.Range("B524").Value = "San Nicolás De Los Garza"
$ws.Range("B527").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B528").Value = "Ayoquezco De Aldama"
$ws.Range("B530").Value = "Coicoyán De Las Flores"
$ws.Range("B532").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B533").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B534").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B537").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B538").Value = "Oaxaca De Juárez"
$ws.Range("B559").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B560").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B587").Value = "Santo Domingo De Morelos"
$ws.Range("B595").Value = "Tataltepec De Valdés"
$ws.Range("B596").Value = "Tlacolula De Matamoros"
$ws.Range("B597").Value = "Villa De Etla"
$ws.Range("B598").Value = "Villa De Tututepec"
$ws.Range("B599").Value = "Villa De Zaachila"
$ws.Range("B600").Value = "Villa Sola De Vega"
$ws.Range("B601").Value = "Zimatlán De Álvarez"
$ws.Range("B620").Value = "Ixcamilpa De Guerrero"
$ws.Range("B622").Value = "Izúcar De Matamoros"
$ws.Range("B626").Value = "Palmar De Bravo"
$ws.Range("B637").Value = "Tecali De Herrera"
$ws.Range("B642").Value = "Tepanco De López"
$ws.Range("B643").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B646").Value = "Tepexi De Rodríguez"
$ws.Range("B660").Value = "Amealco De Bonfil"
$ws.Range("B662").Value = "Cadereyta De Montes"
$ws.Range("B665").Value = "Jalpan De Serra"
$ws.Range("B667").Value = "Pinal De Amoles"
$ws.Range("B669").Value = "San Juan Del Río"
$ws.Range("B689").Value = "San Ciro De Acosta"
$ws.Range("B692").Value = "Santa María Del Río"
$ws.Range("B693").Value = "Soledad De Graciano Sánchez"
$ws.Range("B696").Value = "Villa De Arista"
$ws.Range("B697").Value = "Villa De Arriaga"
$ws.Range("B698").Value = "Villa De La Paz"
$ws.Range("B699").Value = "Villa De Ramos"
$ws.Range("B727").Value = "Jalpa De Méndez"
$ws.Range("B752").Value = "Soto La Marina"
$ws.Range("B759").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B760").Value = "Amaxac De Guerrero"
$ws.Range("B764").Value = "Contla De Juan Cuamatzi"
$ws.Range("B766").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B774").Value = "Tetla De La Solidaridad"
$ws.Range("B785").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B789").Value = "Amatlán De Los Reyes"
$ws.Range("B795").Value = "Castillo De Teayo"
$ws.Range("B802").Value = "Cosamaloapan De Carpio"
$ws.Range("B808").Value = "Hueyapan De Ocampo"
$ws.Range("B809").Value = "Ignacio De La Llave"
$ws.Range("B811").Value = "Ixhuatlán De Madero"
$ws.Range("B812").Value = "Ixhuatlán Del Sureste"
$ws.Range("B817").Value = "Juchique De Ferrer"
$ws.Range("B821").Value = "Lerdo De Tejada"
$ws.Range("B823").Value = "Martínez De La Torre"
$ws.Range("B833").Value = "Paso De Ovejas"
$ws.Range("B836").Value = "Poza Rica De Hidalgo"
$ws.Range("B842").Value = "Sayula De Alemán"
$ws.Range("B843").Value = "Soledad De Doblado"
$ws.Range("B857").Value = "Tlacotepec De Mejía"
$ws.Range("B864").Value = "Vega De Alatorre"
$ws.Range("B876").Value = "Cañitas De Felipe Pescador"
$ws.Range("B878").Value = "Concepción Del Oro"
$ws.Range("B889").Value = "Mezquital Del Oro"
$ws.Range("B893").Value = "Noria De Ángeles"
$ws.Range("B901").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B903").Value = "Villa De Cos"

# Fix tiny floating point value for D732 (Tamaulipas subtotal percentage)
$ws.Range("D732").Value = 0.009442060085836907

# Remove trailing footer/metadata rows (911:915), shrinking used range to A1:D909
$ws.Range("A911:A915").ClearContents()
